# Insert a new row at position 70. This shifts the existing rows 70-81
# down to 71-82, preserving all of their data and formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record's data.
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44637
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112030
$ws.Range("G70").Value = "Poroto granado"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 500
$ws.Range("K70").Value = 18000
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = 19000
$ws.Range("N70").Value = "$/malla 25 kilos"
$ws.Range("O70").Value = "Provincia de Limarí"
$ws.Range("P70").Value = 760
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"

# Make sure the date cell keeps the same number format (date/time) as
# the other rows in column D.
$ws.Range("D70").NumberFormat = $ws.Range("D71").NumberFormat
